# Travelling Salesman Problem - check
# Update the "Fitness" values (column C) on Sheet1 for rows 2-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 4168
    3  = 4332
    4  = 4332
    5  = 4749
    6  = 5059
    7  = 5125
    8  = 5125
    9  = 5125
    10 = 5125
    11 = 5125
    12 = 5125
    13 = 5125
    14 = 5392
    15 = 5392
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
